# Adds "question 10" (SQS) results + the "SQS" tech-label column that
# accompanies them on the other already-partially-filled sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# exp_8 (sheet8.xml): add the "SQS" label column A2:A12, leave rest as-is
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("exp_8")
for ($r = 2; $r -le 12; $r++) {
    $ws8.Range("A$r").Value = "SQS"
}
$ws8.Range("D22").Select()

# ---------------------------------------------------------------------
# exp_9 (sheet9.xml): add the "SQS" label column A2:A22
# ---------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("exp_9")
for ($r = 2; $r -le 22; $r++) {
    $ws9.Range("A$r").Value = "SQS"
}
$ws9.Range("H20").Select()

# ---------------------------------------------------------------------
# exp_10 (sheet10.xml): the actual new question-10 results
# ---------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item("exp_10")

$results = @{
    2  = @(47.416210174560497, 1840.4499567472001, 3599.7326374053901)
    3  = @(2655.94387054443, 4463.19315525201, 6259.5577239990198)
    4  = @(5830.8584690093903, 7679.0393453377901, 9359.6909046173005)
    5  = @(9344.7067737579291, 11072.231920865799, 12817.7707195281)
    6  = @(12758.285999297999, 14580.3457681949, 16452.364683151201)
    7  = @(16401.504755020102, 18241.6293391814, 20045.968770980799)
    8  = @(20047.175407409599, 21819.731088785, 23524.977684020902)
    9  = @(23581.831693649201, 25369.312286376899, 27310.919046401901)
    10 = @(27339.115381240801, 28957.0198013232, 30896.021604537898)
}

for ($r = 2; $r -le 10; $r++) {
    $ws10.Range("A$r").Value = "SQS"
    $vals = $results[$r]
    $ws10.Range("H$r").Value = $vals[0]
    $ws10.Range("I$r").Value = $vals[1]
    $ws10.Range("J$r").Value = $vals[2]
}

# The trailing, empty "staircase" of formatted-but-empty cells (N..R) that
# shows up to the right of the data on every other already-populated
# results sheet (copy/paste artifact from formatting the table). Re-create
# it by copying the matching formatted-but-empty cells from exp_3, which
# already carries the same style indices (s="1" / s="2").
$styleSrc = $wb.Worksheets.Item("exp_3")

$styleSrc.Range("L2").Copy()
$ws10.Range("N2").PasteSpecial(-4122)

for ($r = 3; $r -le 11; $r++) {
    $styleSrc.Range("L4").Copy()
    $ws10.Range("N$r").PasteSpecial(-4122)

    $styleSrc.Range("L2").Copy()
    $ws10.Range("O$r").PasteSpecial(-4122)
    $ws10.Range("P$r").PasteSpecial(-4122)
    $ws10.Range("Q$r").PasteSpecial(-4122)
    $ws10.Range("R$r").PasteSpecial(-4122)
}

$ws10.Range("P19").Select()
